# Updated transition-probability matrix cells after adding more simulated games
# (McNeese_B team-specific Markov transition matrix).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.1755102040816326
$ws.Range("C2").Value = 0.5836734693877551
$ws.Range("J2").Value = 0.00816326530612245
$ws.Range("P2").Value = 0.1102040816326531
$ws.Range("S2").Value = 0.1224489795918367

# Row 3
$ws.Range("B3").Value = 0.01324503311258278
$ws.Range("C3").Value = 0.05298013245033113
$ws.Range("J3").Value = 0.02649006622516556
$ws.Range("P3").Value = 0.7218543046357616
$ws.Range("S3").Value = 0.1854304635761589

# Row 4
$ws.Range("P4").Value = 0.68
$ws.Range("S4").Value = 0.32

# Row 6
$ws.Range("B6").Value = 0.06147540983606557
$ws.Range("E6").Value = 0.004098360655737705
$ws.Range("F6").Value = 0.05737704918032787
$ws.Range("J6").Value = 0.3114754098360656
$ws.Range("O6").Value = 0.01229508196721311
$ws.Range("Q6").Value = 0.110655737704918
$ws.Range("R6").Value = 0.07377049180327869
$ws.Range("S6").Value = 0.3688524590163935

# Row 7
$ws.Range("B7").Value = 0.1212121212121212
$ws.Range("D7").Value = 0.005050505050505051
$ws.Range("F7").Value = 0.03535353535353535
$ws.Range("J7").Value = 0.1464646464646465
$ws.Range("O7").Value = 0.01515151515151515
$ws.Range("Q7").Value = 0.1111111111111111
$ws.Range("R7").Value = 0.101010101010101
$ws.Range("S7").Value = 0.4646464646464646

# Row 8
$ws.Range("B8").Value = 0.06493506493506493
$ws.Range("D8").Value = 0.01082251082251082
$ws.Range("F8").Value = 0.07142857142857142
$ws.Range("J8").Value = 0.1017316017316017
$ws.Range("O8").Value = 0.006493506493506494
$ws.Range("Q8").Value = 0.09740259740259741
$ws.Range("R8").Value = 0.1298701298701299
$ws.Range("S8").Value = 0.5173160173160173

# Row 9
$ws.Range("B9").Value = 0.1160220994475138
$ws.Range("D9").Value = 0.01657458563535912
$ws.Range("E9").Value = 0.005524861878453038
$ws.Range("F9").Value = 0.1104972375690608
$ws.Range("J9").Value = 0.1215469613259668
$ws.Range("O9").Value = 0.01657458563535912
$ws.Range("Q9").Value = 0.08839779005524862
$ws.Range("R9").Value = 0.1215469613259668
$ws.Range("S9").Value = 0.4033149171270718

# Row 10
$ws.Range("B10").Value = 0.1023923444976077
$ws.Range("D10").Value = 0.01626794258373206
$ws.Range("F10").Value = 0.09473684210526316
$ws.Range("J10").Value = 0.1253588516746411
$ws.Range("O10").Value = 0.01052631578947368
$ws.Range("Q10").Value = 0.1464114832535885
$ws.Range("R10").Value = 0.07751196172248803
$ws.Range("S10").Value = 0.4267942583732057

# Row 11
$ws.Range("G11").Value = 0.165625
$ws.Range("J11").Value = 0.08437500000000001
$ws.Range("K11").Value = 0.2375
$ws.Range("L11").Value = 0.490625
$ws.Range("S11").Value = 0.021875

# Row 12
$ws.Range("G12").Value = 0.7409638554216867
$ws.Range("J12").Value = 0.1807228915662651
$ws.Range("K12").Value = 0.01807228915662651
$ws.Range("L12").Value = 0.03012048192771084
$ws.Range("S12").Value = 0.03012048192771084

# Row 13
$ws.Range("F13").Value = 0.02083333333333333
$ws.Range("G13").Value = 0.6458333333333334
$ws.Range("J13").Value = 0.2291666666666667
$ws.Range("S13").Value = 0.1041666666666667

# Row 15
$ws.Range("F15").Value = 0.01639344262295082
$ws.Range("H15").Value = 0.180327868852459
$ws.Range("I15").Value = 0.06010928961748634
$ws.Range("J15").Value = 0.366120218579235
$ws.Range("K15").Value = 0.07650273224043716
$ws.Range("M15").Value = 0.01092896174863388
$ws.Range("O15").Value = 0.06557377049180328
$ws.Range("S15").Value = 0.2240437158469945

# Row 16
$ws.Range("F16").Value = 0.03973509933774835
$ws.Range("H16").Value = 0.1986754966887417
$ws.Range("I16").Value = 0.07947019867549669
$ws.Range("J16").Value = 0.3509933774834437
$ws.Range("K16").Value = 0.1059602649006623
$ws.Range("M16").Value = 0.006622516556291391
$ws.Range("O16").Value = 0.07947019867549669
$ws.Range("S16").Value = 0.1390728476821192

# Row 17
$ws.Range("F17").Value = 0.02290076335877863
$ws.Range("H17").Value = 0.232824427480916
$ws.Range("I17").Value = 0.1145038167938931
$ws.Range("J17").Value = 0.3244274809160305
$ws.Range("K17").Value = 0.0916030534351145
$ws.Range("M17").Value = 0.01908396946564886
$ws.Range("N17").Value = 0.003816793893129771
$ws.Range("O17").Value = 0.05725190839694656
$ws.Range("S17").Value = 0.133587786259542

# Row 18
$ws.Range("F18").Value = 0.0396039603960396
$ws.Range("H18").Value = 0.1732673267326733
$ws.Range("I18").Value = 0.06930693069306931
$ws.Range("J18").Value = 0.3811881188118812
$ws.Range("K18").Value = 0.1633663366336634
$ws.Range("M18").Value = 0.009900990099009901
$ws.Range("N18").Value = 0.004950495049504951
$ws.Range("O18").Value = 0.05445544554455446
$ws.Range("S18").Value = 0.103960396039604

# Row 19
$ws.Range("F19").Value = 0.01805337519623234
$ws.Range("H19").Value = 0.2409733124018838
$ws.Range("I19").Value = 0.08948194662480377
$ws.Range("J19").Value = 0.3194662480376766
$ws.Range("K19").Value = 0.1200941915227629
$ws.Range("M19").Value = 0.03061224489795918
$ws.Range("N19").Value = 0.001569858712715856
$ws.Range("O19").Value = 0.06750392464678179
$ws.Range("S19").Value = 0.1122448979591837
